$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 224). The commit updates this date from 2023-09-06
# (serial 45175) to 2023-09-08 (serial 45177) across all rows.
$ws.Range("C2:C224").Value = 45177
